# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = "29.016.97"
$ws.Range("E2").Value2 = "  -0.90%  "

# Row 3
$ws.Range("D3").Value2 = "1.826.59"
$ws.Range("E3").Value2 = "  -0.92%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.0000"
$ws.Range("E4").Value2 = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "241.67"
$ws.Range("E5").Value2 = "  +0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.6340"
$ws.Range("E6").Value2 = "  -5.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "1.001"
$ws.Range("E7").Value2 = "  +0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "44.66"
$ws.Range("E8").Value2 = "  +6.00%  "

# Row 9
$ws.Range("B9").Value2 = "Cardano"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.2934"
$ws.Range("E9").Value2 = "  -0.16%  "

# Row 10
$ws.Range("B10").Value2 = "Dogecoin"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07340"
$ws.Range("E10").Value2 = "  -1.30%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "22.88"
$ws.Range("E11").Value2 = "  -0.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.07671"
$ws.Range("E12").Value2 = "  -0.64%  "

# Row 13
$ws.Range("D13").Value2 = "1.828.59"
$ws.Range("E13").Value2 = "  -0.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "4.993"
$ws.Range("E14").Value2 = "  -0.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.6653"
$ws.Range("E15").Value2 = "  -0.90%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "82.14"
$ws.Range("E16").Value2 = "  -4.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "6.075"
$ws.Range("E17").Value2 = "  -1.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.000008648"
$ws.Range("E18").Value2 = "  +3.92%  "

# Row 19
$ws.Range("D19").Value2 = "29.019.39"
$ws.Range("E19").Value2 = "  -0.76%  "

# Row 20
$ws.Range("D20").Value2 = "2.082.46"
$ws.Range("E20").Value2 = "  +0.37%  "

# Row 21
$ws.Range("E21").Value2 = "  -0.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "223.99"
$ws.Range("E22").Value2 = "  -1.92%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.001"
$ws.Range("E23").Value2 = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "7.126"
$ws.Range("E24").Value2 = "  -0.52%  "

# Row 25
$ws.Range("E25").Value2 = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "158.14"
$ws.Range("E26").Value2 = "  -1.68%  "

# Row 27
$ws.Range("B27").Value2 = "Cosmos"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "8.455"
$ws.Range("E27").Value2 = "  -2.86%  "

# Row 28
$ws.Range("B28").Value2 = "Stellar"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "0.1376"
$ws.Range("E28").Value2 = "  -1.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.506"
$ws.Range("E30").Value2 = "  -0.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.095"
$ws.Range("E31").Value2 = "  -1.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "4.019"
$ws.Range("E32").Value2 = "  -1.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.201"
$ws.Range("E33").Value2 = "  +0.55%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05308"
$ws.Range("E34").Value2 = "  +0.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.7412"
$ws.Range("E35").Value2 = "  -1.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.832"
$ws.Range("E36").Value2 = "  -2.34%  "

# Row 37
$ws.Range("E37").Value2 = "  +1.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.651"
$ws.Range("E38").Value2 = "  -1.03%  "

# Row 39
$ws.Range("D39").Value2 = "1.288.53"
$ws.Range("E39").Value2 = "  -2.45%  "

# Row 40
$ws.Range("B40").Value2 = "MXToken"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.744"
$ws.Range("E40").Value2 = "  +0.57%  "

# Row 41
$ws.Range("B41").Value2 = "VeChain"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.01782"
$ws.Range("E41").Value2 = "  -1.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "6.347"
$ws.Range("E42").Value2 = "  +6.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.8944"
$ws.Range("E43").Value2 = "  -2.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "1.0000"
$ws.Range("E44").Value2 = "  -0.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "102.67"
$ws.Range("E45").Value2 = "  +0.50%  "

# Row 46
$ws.Range("B46").Value2 = "RocketPoolETH"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value2 = "1.981.12"
$ws.Range("E46").Value2 = "  +0.27%  "

# Row 47
$ws.Range("B47").Value2 = "Mantle"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.5137"
$ws.Range("E47").Value2 = "  -0.51%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "64.20"
$ws.Range("E48").Value2 = "  +0.70%  "

# Row 49
$ws.Range("B49").Value2 = "BabyDogeCoin"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.00000000120"
$ws.Range("E49").Value2 = "  -0.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.730"
$ws.Range("E50").Value2 = "  -2.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.05829"
$ws.Range("E51").Value2 = "  -1.96%  "

